# Weekly fruit/veggie price update:
# insert a new data row at row 204 (shifting the existing rows 204-233 down to
# 205-234) and populate it with a new weekly observation that mirrors the
# values of the (now shifted) row that used to sit at 204, except for the
# date / volume / weighted price / price-per-kg columns which carry the new
# week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 204:233 down one row, to make room for the new
# observation; this also auto-extends the used range down to row 234.
$ws.Rows.Item(204).Insert()

# New row 204 keeps the same Mercado/Region/Categoria/Calidad/etc. as the
# (old) first row of this block, only the date and the three price columns
# change for this week's reading.
$ws.Cells.Item(204, 1).Value = 11
$ws.Cells.Item(204, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(204, 3).Value = "Bíobío"
$ws.Cells.Item(204, 4).Value = 44642
$ws.Cells.Item(204, 5).Value = 8
$ws.Cells.Item(204, 6).Value = 100112009
$ws.Cells.Item(204, 7).Value = "Acelga"
$ws.Cells.Item(204, 8).Value = "Sin especificar"
$ws.Cells.Item(204, 9).Value = "Primera"
$ws.Cells.Item(204, 10).Value = 220
$ws.Cells.Item(204, 11).Value = 600
$ws.Cells.Item(204, 12).Value = 650
$ws.Cells.Item(204, 13).Value = 623
$ws.Cells.Item(204, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(204, 15).Value = "Región de Ñuble"
$ws.Cells.Item(204, 16).Value = 623
$ws.Cells.Item(204, 17).Value = 1
$ws.Cells.Item(204, 18).Value = "Hortaliza"
